$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "_old" suffix -> "_FV2410", "_new" suffix -> "_FV2504"
$oldHeaders = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newHeadersFV2410 = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")

$newSuffixHeaders = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")
$newHeadersFV2504 = @("Segmentname_FV2504","Segmentgruppe_FV2504","Segment_FV2504","Datenelement_FV2504","Segment ID_FV2504","Code_FV2504","Qualifier_FV2504","Beschreibung_FV2504","Bedingungsausdruck_FV2504","Bedingung_FV2504")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2410[$i]
}
# column 11 ("diff") stays the same
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeadersFV2504[$i]
}

# 2. Add a table over the used range, with autofilter + banded rows
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U71"), $false, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row (split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$wb.Save()
